$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "ingrid.matos@mrv.com.br"
$ws.Range("B8").Value = "Excel"
$ws.Range("C8").Value = "Análise de Dados"
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "Teste 3"
